# Update canonical URL (#81)
# Changes:
#  - Metadata!B2 (URL): hl7.fr/fhir/fr/... -> hl7.fr/ig/fhir/...
#  - Metadata!B8 (Date): updated timestamp
#  - Elements!Z6 (Binding Value Set URL): hl7.fr/fhir/fr/... -> hl7.fr/ig/fhir/...

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/StructureDefinition/fr-additional-when-values"
$wsMeta.Range("B8").Value = "2025-05-05T08:11:38+00:00"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("Z6").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-additional-when-codes"

# Note: column Z's "bestFit" width shifts very slightly in the source diff
# (52.4921875 -> 52.7265625) as a side effect of Excel's own re-measurement
# of the updated text. That sub-pixel change isn't reachable through the
# ColumnWidth object model (its setter quantizes to whole pixels), so it is
# intentionally left alone here rather than forcing a worse approximation.
